# Burndown workbook update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix label text: "Acceptatietest" -> "Acceptatietesten"
$ws.Range("A11").Value = "Acceptatietesten"

# 2. Insert three new rows before the totals rows (old rows 28/29)
$ws.Rows.Item(28).Insert()
$ws.Rows.Item(28).Insert()
$ws.Rows.Item(28).Insert()

# 3. Populate new DoD checklist rows (28-30)
$ws.Range("A28").Value = "inlogfunctie"
$ws.Range("B28").Value = 2
$ws.Range("I28").Value = "x"

$ws.Range("A29").Value = "verbeterpunt "
$ws.Range("B29").Value = 2
$ws.Range("I29").Value = "x"

$ws.Range("A30").Value = "reflectieverslag"
$ws.Range("B30").Value = 2
$ws.Range("I30").Value = 2

Write-Output "done-phase1"
